# Auto-generated edit script: refresh '想去人数'/counter values per commit 456a3b4
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("展览")
$ws.Range("F4").Value = 376
$ws.Range("F5").Value = 8313
$ws.Range("G5").Value = 89
$ws.Range("F7").Value = 121
$ws.Range("F12").Value = 639
$ws.Range("F14").Value = 7350
$ws.Range("F16").Value = 7624
$ws.Range("F18").Value = 57549
$ws.Range("F19").Value = 57549
$ws.Range("F20").Value = 4734
$ws.Range("F22").Value = 935
$ws.Range("F23").Value = 496
$ws.Range("F24").Value = 109
$ws.Range("F28").Value = 5290
$ws.Range("F29").Value = 599
$ws.Range("F30").Value = 105
$ws.Range("F32").Value = 903
$ws.Range("F33").Value = 1348
$ws.Range("F34").Value = 1821
$ws.Range("F42").Value = 779
$ws.Range("F43").Value = 256
$ws.Range("F44").Value = 0
$ws.Range("F47").Value = 196

$ws = $wb.Worksheets.Item("演出")
$ws.Range("F2").Value = 197
$ws.Range("F10").Value = 7599
$ws.Range("F24").Value = 31
$ws.Range("F42").Value = 123

$ws = $wb.Worksheets.Item("本地生活")
$ws.Range("F3").Value = 0
$ws.Range("F4").Value = 2360
$ws.Range("F9").Value = 9426
$ws.Range("F10").Value = 1751
$ws.Range("F11").Value = 179
$ws.Range("F16").Value = 2324
$ws.Range("F17").Value = 55
$ws.Range("F18").Value = 491

$ws = $wb.Worksheets.Item("全部类型")
$ws.Range("F3").Value = 2360
$ws.Range("F4").Value = 376
$ws.Range("F5").Value = 8313
$ws.Range("G5").Value = 89
$ws.Range("F7").Value = 179
$ws.Range("F8").Value = 121
$ws.Range("F11").Value = 639
$ws.Range("F13").Value = 7350
$ws.Range("F14").Value = 7624
$ws.Range("F15").Value = 57549
$ws.Range("F16").Value = 197
$ws.Range("F18").Value = 4734
$ws.Range("F20").Value = 935
$ws.Range("F21").Value = 496
$ws.Range("F24").Value = 5290
$ws.Range("F25").Value = 599
$ws.Range("F26").Value = 105
$ws.Range("F27").Value = 903
$ws.Range("F28").Value = 1348
$ws.Range("F30").Value = 491
$ws.Range("F36").Value = 779
$ws.Range("F37").Value = 256
$ws.Range("F45").Value = 196
$ws.Range("F48").Value = 0
